$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "                                         RoboCamp demoyritys 1                                     "
$ws.Range("A3").Value = "                                         Ohjelmistorobotti                                     "

$ws.Range("A3").Select()
